# Update the "想去人数" (F column) counts on the "展览" and "全部类型"
# worksheets, matching the values published for commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Row number (on each sheet) -> new value for column F
$updates = @{
    4  = 349
    6  = 408
    12 = 113
    13 = 1090
    14 = 1452
    15 = 319
    16 = 357
    18 = 93
    20 = 53
    21 = 94
    25 = 1665
    29 = 624
    31 = 3911
    32 = 6
    33 = 455
    34 = 225
    35 = 989
    36 = 96
    39 = 89
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
